$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 28 new rows (119-146) of raw C19 stats data (Datum, Uhrzeit, Summe,
# Summe Aachen, Summe Todesfaelle, Summe genesen, Akute Faelle).
$ws.Range("A119").Value = 44061
$ws.Range("B119").Value = "11:00:00"
$ws.Range("C119").Value = 2248
$ws.Range("D119").Value = 1102
$ws.Range("E119").Value = 101
$ws.Range("F119").Value = 2062
$ws.Range("G119").Value = 85
$ws.Range("A120").Value = 44062
$ws.Range("B120").Value = "11:00:00"
$ws.Range("C120").Value = 2261
$ws.Range("D120").Value = 1106
$ws.Range("E120").Value = 102
$ws.Range("F120").Value = 2074
$ws.Range("G120").Value = 85
$ws.Range("A121").Value = 44063
$ws.Range("B121").Value = "09:30:00"
$ws.Range("C121").Value = 2274
$ws.Range("D121").Value = 1112
$ws.Range("E121").Value = 102
$ws.Range("F121").Value = 2077
$ws.Range("G121").Value = 95
$ws.Range("A122").Value = 44064
$ws.Range("B122").Value = "09:00:00"
$ws.Range("C122").Value = 2281
$ws.Range("D122").Value = 1114
$ws.Range("E122").Value = 102
$ws.Range("F122").Value = 2087
$ws.Range("G122").Value = 92
$ws.Range("A123").Value = 44067
$ws.Range("B123").Value = "10:30:00"
$ws.Range("C123").Value = 2311
$ws.Range("D123").Value = 1126
$ws.Range("E123").Value = 102
$ws.Range("F123").Value = 2120
$ws.Range("G123").Value = 89
$ws.Range("A124").Value = 44068
$ws.Range("B124").Value = "09:30:00"
$ws.Range("C124").Value = 2316
$ws.Range("D124").Value = 1128
$ws.Range("E124").Value = 102
$ws.Range("F124").Value = 2143
$ws.Range("G124").Value = 71
$ws.Range("A125").Value = 44069
$ws.Range("B125").Value = "09:30:00"
$ws.Range("C125").Value = 2324
$ws.Range("D125").Value = 1131
$ws.Range("E125").Value = 102
$ws.Range("F125").Value = 2149
$ws.Range("G125").Value = 73
$ws.Range("A126").Value = 44070
$ws.Range("B126").Value = "09:15:00"
$ws.Range("C126").Value = 2333
$ws.Range("D126").Value = 1136
$ws.Range("E126").Value = 102
$ws.Range("F126").Value = 2156
$ws.Range("G126").Value = 75
$ws.Range("A127").Value = 44071
$ws.Range("B127").Value = "10:30:00"
$ws.Range("C127").Value = 2348
$ws.Range("D127").Value = 1140
$ws.Range("E127").Value = 102
$ws.Range("F127").Value = 2163
$ws.Range("G127").Value = 83
$ws.Range("A128").Value = 44074
$ws.Range("B128").Value = "09:45:00"
$ws.Range("C128").Value = 2353
$ws.Range("D128").Value = 1144
$ws.Range("E128").Value = 102
$ws.Range("F128").Value = 2193
$ws.Range("G128").Value = 58
$ws.Range("A129").Value = 44075
$ws.Range("B129").Value = "09:15:00"
$ws.Range("C129").Value = 2363
$ws.Range("D129").Value = 1149
$ws.Range("E129").Value = 103
$ws.Range("F129").Value = 2203
$ws.Range("G129").Value = 57
$ws.Range("A130").Value = 44076
$ws.Range("B130").Value = "11:30:00"
$ws.Range("C130").Value = 2370
$ws.Range("D130").Value = 1151
$ws.Range("E130").Value = 103
$ws.Range("F130").Value = 2216
$ws.Range("G130").Value = 51
$ws.Range("A131").Value = 44077
$ws.Range("B131").Value = "09:15:00"
$ws.Range("C131").Value = 2378
$ws.Range("D131").Value = 1157
$ws.Range("E131").Value = 103
$ws.Range("F131").Value = 2224
$ws.Range("G131").Value = 51
$ws.Range("A132").Value = 44078
$ws.Range("B132").Value = "09:45:00"
$ws.Range("C132").Value = 2384
$ws.Range("D132").Value = 1161
$ws.Range("E132").Value = 103
$ws.Range("F132").Value = 2235
$ws.Range("G132").Value = 46
$ws.Range("A133").Value = 44081
$ws.Range("B133").Value = "09:45:00"
$ws.Range("C133").Value = 2389
$ws.Range("D133").Value = 1164
$ws.Range("E133").Value = 103
$ws.Range("F133").Value = 2247
$ws.Range("G133").Value = 39
$ws.Range("A134").Value = 44082
$ws.Range("B134").Value = "09:30:00"
$ws.Range("C134").Value = 2393
$ws.Range("D134").Value = 1166
$ws.Range("E134").Value = 103
$ws.Range("F134").Value = 2258
$ws.Range("G134").Value = 32
$ws.Range("A135").Value = 44083
$ws.Range("B135").Value = "09:30:00"
$ws.Range("C135").Value = 2400
$ws.Range("D135").Value = 1171
$ws.Range("E135").Value = 103
$ws.Range("F135").Value = 2260
$ws.Range("G135").Value = 37
$ws.Range("A136").Value = 44084
$ws.Range("B136").Value = "09:30:00"
$ws.Range("C136").Value = 2412
$ws.Range("D136").Value = 1178
$ws.Range("E136").Value = 103
$ws.Range("F136").Value = 2263
$ws.Range("G136").Value = 46
$ws.Range("A137").Value = 44088
$ws.Range("B137").Value = "09:45:00"
$ws.Range("C137").Value = 2429
$ws.Range("D137").Value = 1186
$ws.Range("E137").Value = 103
$ws.Range("F137").Value = 2271
$ws.Range("G137").Value = 55
$ws.Range("A138").Value = 44089
$ws.Range("B138").Value = "08:45:00"
$ws.Range("C138").Value = 2477
$ws.Range("D138").Value = 1210
$ws.Range("E138").Value = 103
$ws.Range("F138").Value = 2288
$ws.Range("G138").Value = 86
$ws.Range("A139").Value = 44090
$ws.Range("B139").Value = "10:45:00"
$ws.Range("C139").Value = 2503
$ws.Range("D139").Value = 1224
$ws.Range("E139").Value = 103
$ws.Range("F139").Value = 2297
$ws.Range("G139").Value = 103
$ws.Range("A140").Value = 44091
$ws.Range("B140").Value = "09:00:00"
$ws.Range("C140").Value = 2522
$ws.Range("D140").Value = 1232
$ws.Range("E140").Value = 103
$ws.Range("F140").Value = 2309
$ws.Range("G140").Value = 110
$ws.Range("A141").Value = 44092
$ws.Range("B141").Value = "09:45:00"
$ws.Range("C141").Value = 2550
$ws.Range("D141").Value = 1248
$ws.Range("E141").Value = 103
$ws.Range("F141").Value = 2326
$ws.Range("G141").Value = 121
$ws.Range("A142").Value = 44095
$ws.Range("B142").Value = "10:00:00"
$ws.Range("C142").Value = 2568
$ws.Range("D142").Value = 1252
$ws.Range("E142").Value = 104
$ws.Range("F142").Value = 2350
$ws.Range("G142").Value = 114
$ws.Range("A143").Value = 44096
$ws.Range("B143").Value = "09:00:00"
$ws.Range("C143").Value = 2584
$ws.Range("D143").Value = 1261
$ws.Range("E143").Value = 104
$ws.Range("F143").Value = 2369
$ws.Range("G143").Value = 111
$ws.Range("A144").Value = 44097
$ws.Range("B144").Value = "09:45:00"
$ws.Range("C144").Value = 2600
$ws.Range("D144").Value = 1268
$ws.Range("E144").Value = 104
$ws.Range("F144").Value = 2387
$ws.Range("G144").Value = 109
$ws.Range("A145").Value = 44098
$ws.Range("B145").Value = "09:00:00"
$ws.Range("C145").Value = 2621
$ws.Range("D145").Value = 1279
$ws.Range("E145").Value = 104
$ws.Range("F145").Value = 2398
$ws.Range("G145").Value = 119
$ws.Range("A146").Value = 44099
$ws.Range("B146").Value = "09:45:00"
$ws.Range("C146").Value = 2643
$ws.Range("D146").Value = 1293
$ws.Range("E146").Value = 104
$ws.Range("F146").Value = 2415
$ws.Range("G146").Value = 124


# Column A (Datum) uses a custom date number format; copy that formatting
# from the last existing row down onto the newly added rows (-4122 ==
# xlPasteFormats, so only styles are pasted, values set above are kept).
$ws.Range("A118").Copy()
$ws.Range("A119:A146").PasteSpecial(-4122)
